$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("menu")

# Row 7: Heart Pan Pizza
$ws.Range("A7").Value = "Heart Pan Pizza"
$ws.Range("B7").Value = 10.1
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = $true

# Row 8: Lasagna
$ws.Range("A8").Value = "Lasagna"
$ws.Range("B8").Value = 69.2
$ws.Range("C8").Value = $true
$ws.Range("D8").Value = $false
